# Weekly driver report update for 2025-04-21
# Applies the updated Critical Minutes / Good Roaming Calculation figures
# for the "Driver Summary" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Row 3: Intel(R) Wi-Fi 6 AX201 160MHz - 22.150.0.3
$ws.Range("C3").Value = 3239
$ws.Range("D3").Value = 76.3

# Row 6: Intel(R) Wi-Fi 6 AX201 160MHz - 22.100.0.3
$ws.Range("D6").Value = 96.59999999999999

# Row 8: Intel(R) Wi-Fi 6 AX201 160MHz - 22.120.0.3
$ws.Range("C8").Value = 742
$ws.Range("D8").Value = 98.3

# Row 10: Totals
$ws.Range("C10").Value = 6512
